# feat: add 2022-Q1 data
#
# - Turns the existing "总计" sheet into the new "2022-Q1" fund-holding
#   sheet (same position/sheetId it already had) and fills it with the
#   Q1-2022 fund detail rows.
# - Adds a fresh "总计" sheet at the end (duplicated from the original
#   "总计" layout) with a new top row summarising 2022-Q1.

$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as text (these tables store numeric-
# looking figures such as "5.67" or fund codes like "002236" as text, not
# numbers) without leaving stray direct formatting behind.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Duplicate "总计" to the end of the workbook, then rename the
#    *original* to "2022-Q1" and the *duplicate* back to "总计". This
#    keeps "总计"'s original sheetId/file slot attached to the new
#    "2022-Q1" sheet (matching how the sheet was actually repurposed) and
#    gives the refreshed "总计" a brand-new sheetId at the end.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Copy($null, $lastSheet)

$q1 = $wb.Worksheets.Item("总计")
$newTotal = $wb.Worksheets.Item("总计 (2)")
$q1.Name = "2022-Q1"
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# 2. Turn "2022-Q1" into a fund-holding table like the other quarters:
#    extend the header row to H1 and replace the old date/count/value
#    rows with the new fund rows.
# ---------------------------------------------------------------------
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2: 002236 - 大成中证360互联网+大数据100指数A
Set-TextValue $q1.Range("B2") "002236"
Set-TextValue $q1.Range("C2") "大成中证360互联网+大数据100指数A"
Set-TextValue $q1.Range("D2") "5.67"
Set-TextValue $q1.Range("E2") "93.32"
Set-TextValue $q1.Range("F2") "1.01"
Set-TextValue $q1.Range("G2") "0.0573"
$q1.Range("H2").Value = 9

# Row 3: 003359 - 大成中证360互联网+大数据100指数C
Set-TextValue $q1.Range("B3") "003359"
Set-TextValue $q1.Range("C3") "大成中证360互联网+大数据100指数C"
Set-TextValue $q1.Range("D3") "4.08"
Set-TextValue $q1.Range("E3") "93.32"
Set-TextValue $q1.Range("F3") "1.01"
Set-TextValue $q1.Range("G3") "0.0412"
$q1.Range("H3").Value = 9

# Row 4 (new row): 001219 - 上投摩根动态多因子策略混合
$q1.Range("A3").Copy()
$q1.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "001219"
Set-TextValue $q1.Range("C4") "上投摩根动态多因子策略混合"
Set-TextValue $q1.Range("D4") "1.16"
Set-TextValue $q1.Range("E4") "92.44"
Set-TextValue $q1.Range("F4") "3.28"
Set-TextValue $q1.Range("G4") "0.0380"
$q1.Range("H4").Value = 3

# ---------------------------------------------------------------------
# 3. Update the refreshed "总计" sheet: insert a new top data row for
#    2022-Q1 and renumber the index column A for the rows pushed down.
# ---------------------------------------------------------------------
$newTotal.Rows("2:2").Insert(-4121)  # xlShiftDown

# The row insert inherits the bold/border formatting of the header row for
# the whole new row; strip that back off of B:D (only column A keeps the
# centered index style, like every other data row in this table).
$newTotal.Range("B2:D2").Style = "Normal"

$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 3
$newTotal.Range("D2").Value = 0.14

$newTotal.Range("A3").Value = 1
$newTotal.Range("A4").Value = 2
$newTotal.Range("A5").Value = 3

# Restore the originally active sheet/selection (untouched by this edit).
$wb.Worksheets.Item("2021-Q2").Activate()
$null = $wb.Worksheets.Item("2021-Q2").Range("A1").Select()
